$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

function Replace-First($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 1) | Out-Null
}

# --- Simple, unambiguous translations -------------------------------------

Replace-All "Observe the effects of surface tension (curved water surface). " "Observez les effets de la tension superficielle (surface de l’eau incurvée). "

Replace-All "Discuss with the other learners why does this effect occur and if they have noticed it before in  daily life." "Discutez avec les autres apprenants de la raison pour laquelle cet effet se produit et s'ils l'ont déjà remarqué dans la vie quotidienne."

Replace-All "Put a layer of water on top of the plates" "Mettre une couche d'eau au-dessus des plaques"

Replace-All "Pour some black pepper or coffee on top of the water. Observe the uniform distribution of the particles" "Verser du poivre noir ou du café sur le dessus de l'eau. Observer la distribution uniforme des particules"

Replace-All "Put a drop of soap on the tip of the sticks" "Mettre une goutte de savon sur la pointe des bâtons"

Replace-All "Touch the water surface with the stick" "Toucher la surface de l'eau avec le bâton"

Replace-All "Observe the spreading of the particles, or, generally, the modification of their distribution" "Observer la propagation des particules ou, en général, la modification de leur distribution"

Replace-All "As coffees are mixtures and their composition can vary, some coffee can react poorly during the experiment." "Comme les cafés sont des mélanges et que leur composition peut varier, certains cafés peuvent réagir mal pendant l’expérience."

Replace-All "After the experiment, the distribution of particles cannot be further modified with the same method." "Après l'expérience, la distribution des particules ne peut plus être modifiée avec la même méthode."

Replace-All "To repeat the experiment, first, clean the plate thoroughly." "Pour répéter l'expérience, d'abord, nettoyer la plaque en profondeur."

# "Invitation to discussion" occurs twice; both get the same French text.
Replace-All "Invitation to discussion" "Invitation à une discussion"

# Longer string first so it doesn't get clobbered by the shorter "Facilitate the discussion" replace below.
Replace-All "Facilitate the discussion: why are the coffee/pepper particles pushed to the plate edge? What is the role of soap?" "Faciliter la discussion : pourquoi les particules de café/poivre sont-elles poussées au bord de la plaque ? Quel est le rôle du savon?"

Replace-All "Suggestion for discussion: surface tension depends on water surface composition " "Suggestion de discussion : la tension de surface dépend de la composition de la surface de l'eau "

# "Try out guesses and share ideas " occurs twice; both get the same French text.
Replace-All "Try out guesses and share ideas " "Essayer de deviner et d'échanger des idées "

Replace-All "Note: pepper or coffee are 'spectators' of the transformation, they serve only to visualize the change in surface tension. The transformation depends only on the soap addition to water." "Remarque : le poivre ou le café sont les « spectateurs » de la transformation, ils servent uniquement à visualiser le changement de tension de surface. La transformation ne dépend que de l’ajout de savon à l’eau."

Replace-All "Important message to deliver: the composition of substances can affect the appearance and properties of objects. The change in composition can manifest itself as a change in the object appearance" "Message important à transmettre : la composition des substances peut affecter l'apparence et les propriétés des objets. Le changement de composition peut se manifester comme un changement dans l'apparence de l'objet"

Replace-All "Experiment solution (part 1)" "Solution d'expérience (partie 1)"

Replace-All "Facilitate the discussion" "Faciliter la discussion"

Replace-All "Suggestion for discussion: surface tension is a surface property" "Suggestion de discussion : la tension de surface dépend de la composition de la surface de l'eau"

Replace-All "Because of their chemical nature, some substances tend to concentrate in specific regions, while some others simply spread randomly" "En raison de leur nature chimique, certaines substances ont tendance à se concentrer dans des régions spécifiques, tandis que d'autres se propagent aléatoirement"

Replace-All "Experiment solution (part 2)" "Solution d'expérience (partie 2)"

# --- "VIDEO PAUSE" multi-segment runs --------------------------------------
# Only the first (of three) "VIDEO PAUSE" occurrences in the document is
# translated; the other two stay in English. Find.Execute with Replace=1
# (wdReplaceOne) over the whole document replaces just the first match
# encountered in document order, which is this one.
Replace-First "VIDEO PAUSE" "PAUSE VIDÉO"

Replace-All "Experiment:" "Expérience :"

Replace-All "Modifying surface tension" "Modification de la tension superficielle"

# Both remaining "Discussion:" occurrences get the same French text.
Replace-All "Discussion:" "Discussion :"

Replace-All "Why do the particles spread?" "Pourquoi les particules se propagent ?"

Replace-All "Where does the soap go?" "Où va le savon ?"
